$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing "Player ID"
# column (and everything else) one column to the right.
$ws.Range("A1").EntireColumn.Insert()

# The header rows (1 & 2) are hidden; writing into a hidden row makes this
# runtime stamp an explicit row height onto it, so temporarily unhide while
# we touch any cells in them and re-hide immediately afterwards.
$row1 = $ws.Rows.Item(1)
$row2 = $ws.Rows.Item(2)
$row20 = $ws.Rows.Item(20)

$row1.Hidden = $false
$row2.Hidden = $false
$row20.Hidden = $false

# New "Match ID" column: header in row 3, a constant match id (27) for every
# player row, and the same match id in the trailing totals row (20).
$ws.Range("A3").Value = "Match ID"
$ws.Range("A4:A19").Value = 27
$ws.Range("A20").Value = 27

$row1.Hidden = $true
$row2.Hidden = $true
$row20.Hidden = $true

# Give the new "Match ID" column the bold, borderless look used elsewhere in
# the sheet for row labels (this creates the same new cell style Excel would
# add on the real edit).
$ws.Range("A3:A19").Font.Bold = $true

# Match the author's resulting selection/viewport.
$ws.Range("A3:A19").Select()
